# Add a new "InvalidLogin" worksheet. Worksheets.Add() inserts the new
# sheet immediately before the currently active sheet (ValidLogin), which
# matches the target layout: InvalidLogin first (sheetId 2), ValidLogin
# second (sheetId 1, keeps its original sheetId).
$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "InvalidLogin"

# Same header row as ValidLogin, but with invalid sample credentials.
$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"

# Match the zoom level used on the ValidLogin sheet.
[void]$newSheet.Range("A1:B2").Select()
$newSheet.Application.ActiveWindow.Zoom = 145

# InvalidLogin is the active/selected tab, with B3 selected.
[void]$newSheet.Range("B3").Select()

# ValidLogin keeps its own selection at A3 (it is no longer the visible tab).
$validSheet = $wb.Worksheets.Item("ValidLogin")
[void]$validSheet.Range("A3").Select()

# Re-activate InvalidLogin so it remains the selected/visible sheet.
[void]$newSheet.Activate()
